# "Latest Specific Workspace and Add On Element"
#
# The Marketo usage report is refreshed for a new workspace
# (".Customer Support New Hire Workspace" / "ACT-SS") and a new
# account ("Pradyumna Sahoo"), and three additional Add-On rows
# (Target Account Management / Predictive Content / Web Personalization)
# are appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Asset Data block (rows 2-7): header + per-item counts ---
$ws.Range("C2").Value = ".Customer Support New Hire Workspace"
$ws.Range("D2").Value = "ACT-SS"
$ws.Range("E2").Value = "'"                 # blank workspace column (text, not 0)

$ws.Range("B3").Value = 7.0                 # Emails
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 7.0
$ws.Range("E3").Value = "'"

$ws.Range("B4").Value = 32.0                # Forms
$ws.Range("C4").Value = 30.0
$ws.Range("D4").Value = 2.0
$ws.Range("E4").Value = "'"

$ws.Range("B5").Value = 71.0                # Landing Pages
$ws.Range("C5").Value = 64.0
$ws.Range("D5").Value = 7.0
$ws.Range("E5").Value = "'"

$ws.Range("B6").Value = 582.0               # Images and Files
$ws.Range("C6").Value = 581.0
$ws.Range("D6").Value = 1.0
$ws.Range("E6").Value = "'"

$ws.Range("B7").Value = 6.0                 # Snippets
$ws.Range("C7").Value = 5.0
$ws.Range("D7").Value = 1.0
$ws.Range("E7").Value = "'"

# --- Campaign Data block (rows 8-14) ---
$ws.Range("C8").Value = ".Customer Support New Hire Workspace"
$ws.Range("D8").Value = "ACT-SS"
$ws.Range("E8").Value = "'"

$ws.Range("B9").Value = 40.0                # All Triggered Campaigns
$ws.Range("C9").Value = "'37"
$ws.Range("D9").Value = "'3"
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = "'"

$ws.Range("B10").Value = 0.0                # Active Triggered Campaigns
$ws.Range("C10").Value = "'0"
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = "'"

# Batch Campaigns - Repeating Schedule (count unchanged, column shifts to F)
$ws.Range("E11").ClearContents()
$ws.Range("F11").Value = "'"

$ws.Range("B12").Value = 84.0               # All Batch Campaigns
$ws.Range("C12").Value = "'82"
$ws.Range("D12").Value = "'2"
$ws.Range("E12").Value = "'"

$ws.Range("B13").Value = 128.0              # All Campaigns
$ws.Range("C13").Value = "'123"
$ws.Range("D13").Value = "'5"
$ws.Range("E13").Value = "'"

$ws.Range("B14").Value = 0.0                # Active Campaigns
$ws.Range("C14").Value = "'0"
$ws.Range("E14").Value = "'"

# --- Database Data block (rows 15-19) ---
$ws.Range("C15").Value = ".Customer Support New Hire Workspace"
$ws.Range("D15").Value = "ACT-SS"

$ws.Range("B16").Value = 0.0                # Segmentations
$ws.Range("C16").Value = 0.0

$ws.Range("B17").Value = 198.0              # Leads
$ws.Range("C17").Value = 0.0
$ws.Range("D17").Value = 0.0

# --- Program Data block (rows 18-19) ---
$ws.Range("C18").Value = ".Customer Support New Hire Workspace"
$ws.Range("D18").Value = "ACT-SS"

$ws.Range("B19").Value = 1.0                # Models
$ws.Range("C19").Value = 1.0

# --- Account summary (rows 20-27) ---
$ws.Range("B20").Value = "'2688"            # Tags
$ws.Range("B21").Value = "'74"              # Integration
$ws.Range("B23").Value = "Pradyumna Sahoo"  # Account Name
$ws.Range("B27").Value = 2.0                # Total WorkSpace

# --- New Add On rows appended at the bottom ---
$ws.Range("A30").Value = "Target Account Management"
$ws.Range("B30").Value = "'True"
$ws.Range("A31").Value = "Predictive Content"
$ws.Range("B31").Value = "'True"
$ws.Range("A32").Value = "Web Personalization"
$ws.Range("B32").Value = "'True"
